$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.248.54"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.342.13"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'585.17"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").Value = "'185.66"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.576"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "'0.182"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "'0.582"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'46.97"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'655.02"
$ws.Range("E13").Value = "  +7.95%  "
$ws.Range("D14").Value = "'8.49"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "3.642.54"
$ws.Range("E15").Value = "  -5.13%  "
$ws.Range("D16").Value = "66.370.82"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "'17.88"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").Value = "3.334.80"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'11.13"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'0.898"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "'17.77"
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("D23").Value = "'5.06"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'100.19"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "'2.80"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'9.67"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").Value = "'32.04"
$ws.Range("E28").Value = "  +5.47%  "
$ws.Range("D29").Value = "'8.56"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'6.83"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "'604.55"
$ws.Range("E31").Value = "  +3.99%  "
$ws.Range("D32").Value = "'3.88"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'11.11"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "3.872.22"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "'56.45"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "'2.76"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'33.00"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0700"
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("D42").Value = "'3.19"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "'0.344"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").Value = "'3.36"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'0.0417"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "'2.56"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").Value = "'2.86"
$ws.Range("E49").Value = "  -17.13%  "
$ws.Range("D50").Value = "'1.33"
$ws.Range("E50").Value = "  +6.78%  "
$ws.Range("D51").Value = "'129.88"
$ws.Range("E51").Value = "  +4.95%  "
